# fix: <aside> páginas de publicaciones
#
# Removes the "5. Casos de error comunes" aside section (its heading plus
# its three bullet rows, along with the blank spacer rows above/below it)
# from the "Git Tag" sheet. Deleting the whole row block shifts every
# subsequent row up and lets Excel compact/renumber the shared-string
# table automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 28-33 hold:
#   28  "5. Casos de error comunes"                                   (header)
#   29  (blank spacer)
#   30  "Qué hacer si se crea un tag con el nombre incorrecto"
#   31  "Problemas de sincronización entre repositorio local y remoto"
#   32  "Conflictos al hacer push de tags"
#   33  (blank spacer)
$ws.Rows("28:33").Delete()

# Land the selection where the author's cursor ended up after the delete.
$ws.Range("A31").Select() | Out-Null
